$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 191, shifting existing rows 191:287 down to 192:288
$ws.Rows("191:191").Insert()

# Populate the newly inserted row 191 with the new data record.
# Columns A,B,C,E,F,G,H,I,R are identical for every record in this sheet
# (same market/region/category), so carry them over unchanged.
$ws.Range("A191").Value = 4
$ws.Range("B191").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C191").Value = "Los Lagos"
$ws.Range("D191").Value = 45001
$ws.Range("E191").Value = 10
$ws.Range("F191").Value = 100112009
$ws.Range("G191").Value = "Acelga"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 35
$ws.Range("K191").Value = 10000
$ws.Range("L191").Value = 10000
$ws.Range("M191").Value = 10000
$ws.Range("N191").Value = "$/docena de atados (12 kilos)"
$ws.Range("O191").Value = "Región de La Araucanía"
$ws.Range("P191").Value = 833
$ws.Range("Q191").Value = 12
$ws.Range("R191").Value = "Hortaliza"
